$d = $word.ActiveDocument

# Locate the paragraph that contains the GaNight astro-map credit line
# (the one still pointing at the 2018 map).
$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*GaNight/2018*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the Jenika Hollana / GaNight credit paragraph"
}

$full = $target.Range
# Replace the paragraph's text (excluding the trailing paragraph mark)
# with the updated credit line that points at the 2022 map.
$r = $d.Range($full.Start, $full.End - 1)
$r.Delete()
$r.InsertAfter(" Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/).")
